# Update the cryptos list with latest price/volume snapshot values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.471.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "'1.896.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'238.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.4902"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").Value = "'0.2925"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("D9").Value = "'0.06690"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("D10").Value = "'1.906.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("D11").Value = "'16.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("D12").Value = "'0.07330"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").Value = "'5.172"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.41%  "
$ws.Range("D14").Value = "'87.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").Value = "'0.6658"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").Value = "'30.435.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'13.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.97%  "
$ws.Range("D18").Value = "'0.000007873"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "'2.158.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("D21").Value = "'5.327"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +13.03%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'193.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "'9.497"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.94%  "
$ws.Range("D26").Value = "'162.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.75%  "
$ws.Range("D27").Value = "'18.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("D28").Value = "'1.939"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.34%  "
$ws.Range("E29").Value = "  +4.88%  "
$ws.Range("D30").Value = "'4.333"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("D31").Value = "'0.09164"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("D32").Value = "'4.062"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.57%  "
$ws.Range("D33").Value = "'0.05162"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("D34").Value = "'0.7410"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.45%  "
$ws.Range("E35").Value = "  +2.17%  "
$ws.Range("D36").Value = "'2.729"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.32%  "
$ws.Range("D37").Value = "'0.01809"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").Value = "'2.688"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").Value = "'0.9236"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.74%  "
$ws.Range("D40").Value = "'2.046"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "'0.4386"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").Value = "'106.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.02%  "
$ws.Range("D43").Value = "'5.919"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.37%  "
$ws.Range("D44").Value = "'0.9937"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "'69.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +20.98%  "
$ws.Range("D46").Value = "'0.1371"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.12%  "
$ws.Range("D47").Value = "'7.576"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("D48").Value = "'9.006"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.86%  "
$ws.Range("D49").Value = "'34.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.46%  "
$ws.Range("D50").Value = "'0.05849"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("D51").Value = "'0.3927"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.15%  "
